$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "In Translation" -> "Handed back: in sync with en-US" everywhere it is
#    used: Overview!E2:F3 (zh-cn / de-de status) and the Status column (C)
#    on the zh-cn and de-de report sheets.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Overview column widths (E, F) widen to fit the longer status text.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# 3. zh-cn report sheet: column widths, handback target/file columns + the
#    handback datetime, plus the new hyperlinks for the Latest Target File
#    column (I).
# ---------------------------------------------------------------------------
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsZh.Range("I2").Value = "34cc9c20-4f37-45de-a9bb-7ff29706a789.md"
$wsZh.Range("J2").Value = "34cc9c20-4f37-45de-a9bb-7ff29706a789.c82444611d0ab2c545bdb763a46541232a8b4edf.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-16 06:21:13"

$wsZh.Range("I3").Value = "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md"
$wsZh.Range("J3").Value = "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.cd692510424258e6612deda64b4baa7ca38eea36.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-16 06:21:13"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba1008df9c021af4f8e91a77c1d6b75f263c031/e2e/34cc9c20-4f37-45de-a9bb-7ff29706a789.md", "", "", "34cc9c20-4f37-45de-a9bb-7ff29706a789.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba1008df9c021af4f8e91a77c1d6b75f263c031/e2e/c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md", "", "", "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md")

$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# 4. de-de report sheet: same shape of edits, with its own handback
#    datetime value.
# ---------------------------------------------------------------------------
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

$wsDe.Range("I2").Value = "34cc9c20-4f37-45de-a9bb-7ff29706a789.md"
$wsDe.Range("J2").Value = "34cc9c20-4f37-45de-a9bb-7ff29706a789.c82444611d0ab2c545bdb763a46541232a8b4edf.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-16 06:21:20"

$wsDe.Range("I3").Value = "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md"
$wsDe.Range("J3").Value = "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.cd692510424258e6612deda64b4baa7ca38eea36.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-16 06:21:20"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba1008df9c021af4f8e91a77c1d6b75f263c031/e2e/34cc9c20-4f37-45de-a9bb-7ff29706a789.md", "", "", "34cc9c20-4f37-45de-a9bb-7ff29706a789.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba1008df9c021af4f8e91a77c1d6b75f263c031/e2e/c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md", "", "", "c42a9e4e-d16b-4ab4-bce0-bdf486f64fe2.md")

$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = 15570276

Write-Host "Handback report generated."
